$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = 0.2300884955752212
$ws.Range("C2").Value2 = 0.4690265486725664
$ws.Range("J2").Value2 = 0.004424778761061947
$ws.Range("P2").Value2 = 0.1814159292035398
$ws.Range("S2").Value2 = 0.1150442477876106
$ws.Range("B3").Value2 = 0.009174311926605505
$ws.Range("C3").Value2 = 0.02752293577981652
$ws.Range("J3").Value2 = 0.02752293577981652
$ws.Range("P3").Value2 = 0.8165137614678899
$ws.Range("S3").Value2 = 0.1192660550458716
$ws.Range("P4").Value2 = 0.5416666666666666
$ws.Range("S4").Value2 = 0.4583333333333333
$ws.Range("B6").Value2 = 0.05851063829787234
$ws.Range("D6").Value2 = 0.01063829787234043
$ws.Range("F6").Value2 = 0.07446808510638298
$ws.Range("J6").Value2 = 0.2180851063829787
$ws.Range("O6").Value2 = 0.02659574468085106
$ws.Range("Q6").Value2 = 0.148936170212766
$ws.Range("R6").Value2 = 0.1170212765957447
$ws.Range("S6").Value2 = 0.3457446808510639
$ws.Range("B7").Value2 = 0.1337579617834395
$ws.Range("D7").Value2 = 0.006369426751592357
$ws.Range("E7").Value2 = 0.006369426751592357
$ws.Range("F7").Value2 = 0.05732484076433121
$ws.Range("J7").Value2 = 0.1401273885350318
$ws.Range("O7").Value2 = 0.01273885350318471
$ws.Range("Q7").Value2 = 0.1528662420382166
$ws.Range("R7").Value2 = 0.08280254777070063
$ws.Range("S7").Value2 = 0.4076433121019108
$ws.Range("B8").Value2 = 0.08232445520581114
$ws.Range("D8").Value2 = 0.01452784503631961
$ws.Range("E8").Value2 = 0.002421307506053269
$ws.Range("F8").Value2 = 0.06053268765133172
$ws.Range("J8").Value2 = 0.08716707021791767
$ws.Range("O8").Value2 = 0.007263922518159807
$ws.Range("Q8").Value2 = 0.198547215496368
$ws.Range("R8").Value2 = 0.07263922518159806
$ws.Range("S8").Value2 = 0.4745762711864407
$ws.Range("B9").Value2 = 0.09722222222222222
$ws.Range("D9").Value2 = 0.01388888888888889
$ws.Range("F9").Value2 = 0.04629629629629629
$ws.Range("J9").Value2 = 0.1203703703703704
$ws.Range("O9").Value2 = 0.009259259259259259
$ws.Range("Q9").Value2 = 0.1527777777777778
$ws.Range("R9").Value2 = 0.07407407407407407
$ws.Range("S9").Value2 = 0.4861111111111111
$ws.Range("B10").Value2 = 0.07721280602636535
$ws.Range("D10").Value2 = 0.01035781544256121
$ws.Range("E10").Value2 = 0.0009416195856873823
$ws.Range("F10").Value2 = 0.07909604519774012
$ws.Range("J10").Value2 = 0.09981167608286252
$ws.Range("O10").Value2 = 0.01318267419962335
$ws.Range("Q10").Value2 = 0.2231638418079096
$ws.Range("R10").Value2 = 0.07250470809792843
$ws.Range("S10").Value2 = 0.423728813559322
$ws.Range("G11").Value2 = 0.1483050847457627
$ws.Range("J11").Value2 = 0.08050847457627118
$ws.Range("K11").Value2 = 0.2033898305084746
$ws.Range("L11").Value2 = 0.559322033898305
$ws.Range("S11").Value2 = 0.008474576271186441
$ws.Range("G12").Value2 = 0.7481481481481481
$ws.Range("J12").Value2 = 0.1703703703703704
$ws.Range("K12").Value2 = 0.007407407407407408
$ws.Range("L12").Value2 = 0.02962962962962963
$ws.Range("S12").Value2 = 0.04444444444444445
$ws.Range("G13").Value2 = 0.6764705882352942
$ws.Range("J13").Value2 = 0.2647058823529412
$ws.Range("S13").Value2 = 0.05882352941176471
$ws.Range("F15").Value2 = 0.01515151515151515
$ws.Range("H15").Value2 = 0.1717171717171717
$ws.Range("I15").Value2 = 0.101010101010101
$ws.Range("J15").Value2 = 0.3686868686868687
$ws.Range("K15").Value2 = 0.0505050505050505
$ws.Range("M15").Value2 = 0.0101010101010101
$ws.Range("N15").Value2 = 0.005050505050505051
$ws.Range("O15").Value2 = 0.05555555555555555
$ws.Range("S15").Value2 = 0.2222222222222222
$ws.Range("H16").Value2 = 0.1690140845070423
$ws.Range("I16").Value2 = 0.09859154929577464
$ws.Range("J16").Value2 = 0.5140845070422535
$ws.Range("K16").Value2 = 0.07042253521126761
$ws.Range("M16").Value2 = 0.01408450704225352
$ws.Range("O16").Value2 = 0.04225352112676056
$ws.Range("S16").Value2 = 0.09154929577464789
$ws.Range("F17").Value2 = 0.01240694789081886
$ws.Range("H17").Value2 = 0.1861042183622829
$ws.Range("I17").Value2 = 0.1042183622828784
$ws.Range("J17").Value2 = 0.3970223325062035
$ws.Range("K17").Value2 = 0.08436724565756824
$ws.Range("M17").Value2 = 0.02233250620347394
$ws.Range("N17").Value2 = 0.004962779156327543
$ws.Range("O17").Value2 = 0.08933002481389578
$ws.Range("S17").Value2 = 0.09925558312655088
$ws.Range("F18").Value2 = 0.02547770700636943
$ws.Range("H18").Value2 = 0.2229299363057325
$ws.Range("I18").Value2 = 0.1273885350318471
$ws.Range("J18").Value2 = 0.3821656050955414
$ws.Range("K18").Value2 = 0.08280254777070063
$ws.Range("M18").Value2 = 0.01273885350318471
$ws.Range("O18").Value2 = 0.08917197452229299
$ws.Range("S18").Value2 = 0.05732484076433121
$ws.Range("F19").Value2 = 0.01492537313432836
$ws.Range("H19").Value2 = 0.2159789288849868
$ws.Range("I19").Value2 = 0.1088674275680421
$ws.Range("J19").Value2 = 0.3792800702370501
$ws.Range("K19").Value2 = 0.1027216856892011
$ws.Range("M19").Value2 = 0.01668129938542581
$ws.Range("N19").Value2 = 0.002633889376646181
$ws.Range("O19").Value2 = 0.06672519754170325
$ws.Range("S19").Value2 = 0.09218612818261633
